$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-10 Saturday" "2024-08-11 Sunday"
Replace-Text "752÷3=250, 2" "287÷8=35, 7"
Replace-Text "346÷4=86, 2" "448÷9=49, 7"
Replace-Text "722÷8=90, 2" "674÷5=134, 4"
Replace-Text "587÷8=73, 3" "432÷6=72, 0"
Replace-Text "832÷9=92, 4" "540÷4=135, 0"
Replace-Text "159÷2=79, 1" "249÷8=31, 1"
Replace-Text "249÷9=27, 6" "882÷2=441, 0"
Replace-Text "865÷9=96, 1" "210÷2=105, 0"
Replace-Text "629÷8=78, 5" "745÷6=124, 1"
Replace-Text "880÷4=220, 0" "305÷2=152, 1"
Replace-Text "418÷3=139, 1" "123÷9=13, 6"
Replace-Text "978÷5=195, 3" "776÷3=258, 2"
Replace-Text "375÷7=53, 4" "903÷9=100, 3"
Replace-Text "633÷7=90, 3" "289÷8=36, 1"
Replace-Text "706÷4=176, 2" "670÷9=74, 4"
Replace-Text "710÷8=88, 6" "411÷7=58, 5"
Replace-Text "494÷5=98, 4" "540÷8=67, 4"
Replace-Text "509÷5=101, 4" "413÷4=103, 1"
Replace-Text "385÷5=77, 0" "836÷2=418, 0"
Replace-Text "520÷7=74, 2" "608÷6=101, 2"
Replace-Text "628÷4=157, 0" "155÷2=77, 1"
Replace-Text "494÷6=82, 2" "365÷4=91, 1"
Replace-Text "844÷9=93, 7" "815÷4=203, 3"
Replace-Text "429÷5=85, 4" "606÷5=121, 1"
Replace-Text "250÷9=27, 7" "949÷4=237, 1"
